# Repull data, push all data, mean calculation
# Update the dSF (column F) values for the affected rows to reflect the
# re-pulled / recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = -1
    12 = -1
    14 = -7
    15 = -4
    17 = -4
    18 = 3
    21 = -6
    23 = 0
    24 = 3
    26 = 0
    31 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
